# Table 2 - SES Korsoer
# Add all-cause monthly mort for CPH and all other cities
# Faceted plots with differering colors and linetypes
#
# - Replace "Attack rate (%)" / "Mortality rate (%)" numeric columns with
#   text "ratio (95% CI)" columns (Attack ratio / Mortality ratio).
# - Re-style the whole table with Times New Roman, widen F/G columns,
#   taller rows, right-aligned+wrapped ratio cells, bottom border moved to
#   the "Low" row (now the last data row before Totals).
# - Row 6 / column C ("Totals" / "Mean house value") becomes "NA" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Column widths for the new/wider ratio columns
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 12.77734375
$ws.Columns.Item(7).ColumnWidth = 13.88671875

# ---------------------------------------------------------------------
# 2. Header row text (F1, G1 change; rest stay the same)
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Attack ratio (95% CI)"
$ws.Range("G1").Value = "Mortality ratio (95% CI)"

# ---------------------------------------------------------------------
# 3. Data values: F/G become text "x.xx`n(lo - hi)" ratio cells,
#    C6 becomes "NA"
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "0.08`n(0.06 - 0.10)"
$ws.Range("G2").Value = "0.05`n(0.03 - 0.06)"

$ws.Range("F3").Value = "0.14`n(0.11 - 0.17)"
$ws.Range("G3").Value = "0.10`n(0.07 - 0.12)"

$ws.Range("F4").Value = "0.17`n(0.13 - 0.21)"
$ws.Range("G4").Value = "0.12`n(0.08 - 0.16)"

$ws.Range("F5").Value = "0.21`n(0.16 - 0.25)"
$ws.Range("G5").Value = "0.14`n(0.10 - 0.18)"

$ws.Range("F6").Value = "0.13`n(0.12 - 0.15)"
$ws.Range("G6").Value = "0.09`n(0.08 - 0.10)"

$ws.Range("C6").Value = "NA"

# ---------------------------------------------------------------------
# 4. Fonts: whole table moves from Calibri to Times New Roman (bold stays
#    bold on the header row, regular elsewhere)
# ---------------------------------------------------------------------
$used = $ws.Range("A1:H6")
$used.Font.Name = "Times New Roman"
$used.Font.Family = 1

# ---------------------------------------------------------------------
# 5. Number formats / alignment for the ratio (F,G) and CFR (H) columns
# ---------------------------------------------------------------------
$ratios = $ws.Range("F2:G6")
$ratios.NumberFormat = "0.00"
$ratios.HorizontalAlignment = -4152   # xlRight
$ratios.WrapText = $true

$ws.Range("H2:H6").NumberFormat = "0%"

# ---------------------------------------------------------------------
# 6. Row heights
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 28.2
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30

# ---------------------------------------------------------------------
# 7. Borders: the thick separator used to sit under row 2 ("High") /
#    above row 3 because row 2 was the header; now the header is row 1
#    and the medium rule sits under row 5 ("Low"), the last detail row
#    before "Totals". Clear the old mid-table border and the stray top
#    border that used to sit above the "Low" row, then draw the medium
#    rule under row 5.
# ---------------------------------------------------------------------
$ws.Range("A2:H4").Borders.LineStyle = -4142   # xlLineStyleNone
$ws.Range("A5:H5").Borders.LineStyle = -4142   # xlLineStyleNone
$ws.Range("A6:H6").Borders.LineStyle = -4142   # xlLineStyleNone

$ws.Range("A5:H5").Borders.Item(9).LineStyle = 1       # xlContinuous
$ws.Range("A5:H5").Borders.Item(9).Weight = -4138      # xlMedium

Write-Output "applied table restyle"
